$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting (values like "1.003" or
# "0.2570" must not be auto-coerced into numbers by COM automation).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '25.912.45'
$ws.Range('E2').Value = '  -0.92%  '

# Row 3
$ws.Range('D3').Value = '1.640.34'
$ws.Range('E3').Value = '  -0.64%  '

# Row 4
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.19%  '

# Row 5
$ws.Range('D5').Value = '216.01'
$ws.Range('E5').Value = '  +0.22%  '

# Row 6
$ws.Range('D6').Value = '0.5039'
$ws.Range('E6').Value = '  -1.59%  '

# Row 7
$ws.Range('E7').Value = '  -0.09%  '

# Row 8
$ws.Range('D8').Value = '0.2570'
$ws.Range('E8').Value = '  -0.81%  '

# Row 9
$ws.Range('D9').Value = '0.06401'
$ws.Range('E9').Value = '  -0.40%  '

# Row 10
$ws.Range('D10').Value = '19.72'
$ws.Range('E10').Value = '  -0.96%  '

# Row 11
$ws.Range('D11').Value = '0.07737'
$ws.Range('E11').Value = '  -0.65%  '

# Row 12
$ws.Range('D12').Value = '4.271'
$ws.Range('E12').Value = '  -0.30%  '

# Row 13
$ws.Range('D13').Value = '1.644.73'
$ws.Range('E13').Value = '  -0.36%  '

# Row 14
$ws.Range('D14').Value = '1.863.82'
$ws.Range('E14').Value = '  -0.71%  '

# Row 15
$ws.Range('D15').Value = '0.5460'
$ws.Range('E15').Value = '  -0.94%  '

# Row 16
$ws.Range('D16').Value = '0.0₅7925'
$ws.Range('E16').Value = '  -0.97%  '

# Row 17
$ws.Range('D17').Value = '64.46'
$ws.Range('E17').Value = '  +0.66%  '

# Row 18
$ws.Range('D18').Value = '25.930.45'
$ws.Range('E18').Value = '  -0.84%  '

# Row 19
$ws.Range('E19').Value = '  -0.06%  '

# Row 20
$ws.Range('D20').Value = '203.66'
$ws.Range('E20').Value = '  -3.02%  '

# Row 21
$ws.Range('D21').Value = '4.402'
$ws.Range('E21').Value = '  +0.06%  '

# Row 22
$ws.Range('D22').Value = '9.935'
$ws.Range('E22').Value = '  -1.14%  '

# Row 23
$ws.Range('D23').Value = '5.991'
$ws.Range('E23').Value = '  -0.75%  '

# Row 24
$ws.Range('D24').Value = '1.006'
$ws.Range('E24').Value = '  +0.08%  '

# Row 25
$ws.Range('D25').Value = '1.936'
$ws.Range('E25').Value = '  +10.84%  '

# Row 26
$ws.Range('D26').Value = '141.48'
$ws.Range('E26').Value = '  -1.77%  '

# Row 27
$ws.Range('D27').Value = '0.1138'
$ws.Range('E27').Value = '  -3.13%  '

# Row 28
$ws.Range('D28').Value = '15.73'
$ws.Range('E28').Value = '  -0.36%  '

# Row 29
$ws.Range('D29').Value = '6.747'
$ws.Range('E29').Value = '  -3.07%  '

# Row 30
$ws.Range('D30').Value = '1.247'
$ws.Range('E30').Value = '  +0.49%  '

# Row 31
$ws.Range('D31').Value = '0.04949'
$ws.Range('E31').Value = '  -3.55%  '

# Row 32
$ws.Range('D32').Value = '3.284'
$ws.Range('E32').Value = '  -1.98%  '

# Row 33
$ws.Range('D33').Value = '3.191'
$ws.Range('E33').Value = '  -0.67%  '

# Row 34
$ws.Range('D34').Value = '1.548'
$ws.Range('E34').Value = '  -0.63%  '

# Row 35
$ws.Range('D35').Value = '2.381'
$ws.Range('E35').Value = '  +1.23%  '

# Row 36
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '0.8972'
$ws.Range('E36').Value = '  -2.88%  '

# Row 37
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').Value = '2.635'
$ws.Range('E37').Value = '  -3.73%  '

# Row 38
$ws.Range('D38').Value = '1.167.31'
$ws.Range('E38').Value = '  +0.44%  '

# Row 39
$ws.Range('D39').Value = '0.5617'
$ws.Range('E39').Value = '  -1.44%  '

# Row 40
$ws.Range('D40').Value = '0.01567'
$ws.Range('E40').Value = '  -1.11%  '

# Row 41
$ws.Range('D41').Value = '1.006'
$ws.Range('E41').Value = '  +0.10%  '

# Row 42
$ws.Range('D42').Value = '5.739'
$ws.Range('E42').Value = '  +1.53%  '

# Row 43
$ws.Range('D43').Value = '0.8103'
$ws.Range('E43').Value = '  -1.49%  '

# Row 44
$ws.Range('D44').Value = '100.01'
$ws.Range('E44').Value = '  -0.10%  '

# Row 45
$ws.Range('D45').Value = '1.775.69'
$ws.Range('E45').Value = '  -0.71%  '

# Row 46
$ws.Range('D46').Value = '0.0₈116'
$ws.Range('E46').Value = '  +0.11%  '

# Row 47
$ws.Range('D47').Value = '0.4525'
$ws.Range('E47').Value = '  -0.49%  '

# Row 48
$ws.Range('E48').Value = '  +0.02%  '

# Row 49
$ws.Range('D49').Value = '55.16'
$ws.Range('E49').Value = '  -0.49%  '

# Row 50
$ws.Range('D50').Value = '0.05056'
$ws.Range('E50').Value = '  -0.42%  '

# Row 51
$ws.Range('D51').Value = '1.003'
$ws.Range('E51').Value = '  -0.30%  '
